# Fruta / hortaliza, semanal
# Insert two new weekly price-report rows for "Membrillo" (Vega Modelo de
# Temuco) at the top of the existing block (old rows 171-184), pushing the
# existing rows down by two (to 173-186).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert two blank rows above the current row 171.
$ws.Rows("171:172").Insert()

# New row 171
$ws.Range("A171").Value = 10
$ws.Range("B171").Value = "Vega Modelo de Temuco"
$ws.Range("C171").Value = "La Araucanía"
$ws.Range("D171").Value = 44746
$ws.Range("E171").Value = 9
$ws.Range("F171").Value = "Fruta"
$ws.Range("G171").Value = 100104
$ws.Range("H171").Value = "Frutos de pepita"
$ws.Range("I171").Value = 100104003
$ws.Range("J171").Value = "Membrillo"
$ws.Range("K171").Value = "Champion"
$ws.Range("L171").Value = "Primera"
$ws.Range("M171").Value = 150
$ws.Range("N171").Value = 10000
$ws.Range("O171").Value = 10000
$ws.Range("P171").Value = 10000
$ws.Range("Q171").Value = "$/bandeja 18 kilos granel"
$ws.Range("R171").Value = "Región de O'Higgins"
$ws.Range("S171").Value = 556
$ws.Range("T171").Value = 18

# New row 172
$ws.Range("A172").Value = 10
$ws.Range("B172").Value = "Vega Modelo de Temuco"
$ws.Range("C172").Value = "La Araucanía"
$ws.Range("D172").Value = 44746
$ws.Range("E172").Value = 9
$ws.Range("F172").Value = "Fruta"
$ws.Range("G172").Value = 100104
$ws.Range("H172").Value = "Frutos de pepita"
$ws.Range("I172").Value = 100104003
$ws.Range("J172").Value = "Membrillo"
$ws.Range("K172").Value = "Champion"
$ws.Range("L172").Value = "Primera"
$ws.Range("M172").Value = 2
$ws.Range("N172").Value = 200000
$ws.Range("O172").Value = 200000
$ws.Range("P172").Value = 200000
$ws.Range("Q172").Value = "$/bins (450 kilos)"
$ws.Range("R172").Value = "Región de O'Higgins"
$ws.Range("S172").Value = 444
$ws.Range("T172").Value = 450
